$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Merges first: merging a range *after* it already has borders applied makes
# this host redistribute the border into separate top/middle/bottom pieces
# (like real Excel's "merge adjusts inner borders" behaviour). The edit we
# are reproducing keeps one uniform box border (borderId 1) on every cell, so
# merge the ranges before any border/font/alignment formatting is applied.
# ---------------------------------------------------------------------------

$ws.Range("B3:B5").Merge()
$ws.Range("A3:A5").Merge()
$ws.Range("B6").Merge()
$ws.Range("A6").Merge()

# ---------------------------------------------------------------------------
# New cell formatting (font sz=14, thin black border everywhere) with four
# alignment variants, matching the four new cellXfs added by the edit:
#   - center/center + wrap  -> A3, A6   (short "model" codes)
#   - center/center         -> B3,D3,D4,D5,B6,D6 (amount/rodzaj text)
#   - right/center          -> C3,C4,C5,C6 (ilosc numbers)
#   - no alignment           -> E3,A4,B4,E4,A5,B5,E5,E6 (empty cells)
# Apply ALL formatting for ALL cells first, then fill in values afterwards -
# interleaving value writes with style writes on different cells can cause
# this host's style de-duplication to miss reusing an identical xf.
# ---------------------------------------------------------------------------

$xlCenter = -4108
$xlRight = -4152

$centerWrapCells = @("A3", "A6")
$centerCells = @("B3", "D3", "D4", "D5", "B6", "D6")
$rightCells = @("C3", "C4", "C5", "C6")
$plainCells = @("E3", "A4", "B4", "E4", "A5", "B5", "E5", "E6")

foreach ($c in $centerWrapCells) {
    $r = $ws.Range($c)
    $r.Font.Size = 14
    $r.Borders.LineStyle = 1
    $r.Borders.Color = 0
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    $r.WrapText = $true
}

foreach ($c in $centerCells) {
    $r = $ws.Range($c)
    $r.Font.Size = 14
    $r.Borders.LineStyle = 1
    $r.Borders.Color = 0
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
}

foreach ($c in $rightCells) {
    $r = $ws.Range($c)
    $r.Font.Size = 14
    $r.Borders.LineStyle = 1
    $r.Borders.Color = 0
    $r.HorizontalAlignment = $xlRight
    $r.VerticalAlignment = $xlCenter
}

foreach ($c in $plainCells) {
    $r = $ws.Range($c)
    $r.Font.Size = 14
    $r.Borders.LineStyle = 1
    $r.Borders.Color = 0
}

# ---------------------------------------------------------------------------
# Values - written only after every cell already carries its final style.
# Leading "'" forces a pure-numeric-looking entry ("27", "50") to stay text
# instead of being auto-converted to a number, matching the inlineStr cells
# from the edit.
# ---------------------------------------------------------------------------

$ws.Range("A3").Value = "D2"
$ws.Range("B3").Value = "'27"
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = "czarny"

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "czarny"

$ws.Range("C5").Value = 12
$ws.Range("D5").Value = "bialy"

$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "'50"
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = "styropian"
